$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$borderColor = 13882323   # RGB(211,211,211) -> 0xD3D3D3, stored BGR in VBA as 0xBBGGRR
$headerFill  = 12180223   # RGB(255,218,185) -> 0xFFDAB9 new header shading color

$rowCount = $t.Rows.Count
$colCount = $t.Columns.Count

for ($r = 1; $r -le $rowCount; $r++) {
  for ($c = 1; $c -le $colCount; $c++) {
    $cell = $t.Cell($r, $c)

    # Add a thin light-grey border on all four sides of every cell.
    $cell.Borders.Item(-1).LineStyle = 1
    $cell.Borders.Item(-1).LineWidth = 2
    $cell.Borders.Item(-1).Color = $borderColor

    $cell.Borders.Item(-2).LineStyle = 1
    $cell.Borders.Item(-2).LineWidth = 2
    $cell.Borders.Item(-2).Color = $borderColor

    $cell.Borders.Item(-3).LineStyle = 1
    $cell.Borders.Item(-3).LineWidth = 2
    $cell.Borders.Item(-3).Color = $borderColor

    $cell.Borders.Item(-4).LineStyle = 1
    $cell.Borders.Item(-4).LineWidth = 2
    $cell.Borders.Item(-4).Color = $borderColor

    if ($r -eq 1) {
      # Header row: re-tint the shading from peach (FFEFD5) to a deeper peach (FFDAB9).
      $cell.Shading.BackgroundPatternColor = $headerFill
    } else {
      # Data rows: every column except the first (label) column is recentered.
      if ($c -ne 1) {
        $cell.Range.ParagraphFormat.Alignment = 1
      }
    }
  }
}
